$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(22)
$r = $p.Range
$h = $d.Hyperlinks.Add($r, "https://www.tinkercad.com/things/8CotskOeSCb")

$p2 = $d.Paragraphs.Item(22)
$p2.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:hyperlink r:id="rId9" w:history="1"><w:r><w:rPr><w:rStyle w:val="SomeUnknownStyleXYZ"/></w:rPr><w:t>https://www.tinkercad.com/things/8CotskOeSCb</w:t></w:r></w:hyperlink></w:p>', "Replace")
Write-Host "done"
